$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain text even though the new values look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.669.64"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "3.252.69"
$ws.Range("E3").Value = "  +6.55%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "581.42"
$ws.Range("E5").Value = "  +4.23%  "
$ws.Range("D6").Value = "153.96"
$ws.Range("E6").Value = "  +8.48%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "3.245.00"
$ws.Range("E8").Value = "  +6.70%  "
$ws.Range("D9").Value = "0.515"
$ws.Range("E9").Value = "  +4.34%  "
$ws.Range("D10").Value = "7.11"
$ws.Range("E10").Value = "  +9.08%  "
$ws.Range("D11").Value = "0.165"
$ws.Range("E11").Value = "  +5.24%  "
$ws.Range("D12").Value = "0.490"
$ws.Range("E12").Value = "  +4.04%  "
$ws.Range("D13").Value = "37.86"
$ws.Range("E13").Value = "  +3.08%  "
$ws.Range("D14").Value = "0.0000235"
$ws.Range("E14").Value = "  +5.13%  "
$ws.Range("D15").Value = "3.773.25"
$ws.Range("E15").Value = "  +6.33%  "
$ws.Range("D16").Value = "558.20"
$ws.Range("E16").Value = "  +12.22%  "
$ws.Range("D17").Value = "66.646.15"
$ws.Range("E17").Value = "  +3.44%  "
$ws.Range("D18").Value = "3.246.47"
$ws.Range("E18").Value = "  +6.09%  "
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").Value = "7.11"
$ws.Range("E20").Value = "  +5.13%  "
$ws.Range("D21").Value = "14.42"
$ws.Range("E21").Value = "  +4.13%  "
$ws.Range("D22").Value = "0.743"
$ws.Range("E22").Value = "  +6.98%  "
$ws.Range("D23").Value = "7.78"
$ws.Range("E23").Value = "  +7.40%  "
$ws.Range("D24").Value = "13.65"
$ws.Range("E24").Value = "  +6.23%  "
$ws.Range("D25").Value = "81.91"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "9.31"
$ws.Range("E27").Value = "  +18.67%  "
$ws.Range("D28").Value = "2.97"
$ws.Range("E28").Value = "  +6.81%  "
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("D30").Value = "27.86"
$ws.Range("E30").Value = "  +6.11%  "
$ws.Range("D31").Value = "2.75"
$ws.Range("E31").Value = "  +3.66%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("E33").Value = "  +5.02%  "
$ws.Range("D34").Value = "563.23"
$ws.Range("E34").Value = "  +9.17%  "
$ws.Range("D35").Value = "5.73"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "6.39"
$ws.Range("E36").Value = "  +5.92%  "
$ws.Range("D37").Value = "0.0459"
$ws.Range("E37").Value = "  +12.51%  "
$ws.Range("D38").Value = "55.56"
$ws.Range("E38").Value = "  +4.82%  "
$ws.Range("D39").Value = "0.0864"
$ws.Range("E39").Value = "  +6.96%  "
$ws.Range("E40").Value = "  +7.06%  "
$ws.Range("D41").Value = "3.06"
$ws.Range("E41").Value = "  +14.57%  "
$ws.Range("D42").Value = "3.166.96"
$ws.Range("E42").Value = "  +7.91%  "
$ws.Range("D43").Value = "8.62"
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("D44").Value = "0.275"
$ws.Range("E44").Value = "  +10.83%  "
$ws.Range("E45").Value = "  +7.40%  "
$ws.Range("D46").Value = "26.61"
$ws.Range("E46").Value = "  +4.69%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "0.0₃0557"
$ws.Range("E48").Value = "  +2.62%  "
$ws.Range("D49").Value = "125.68"
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").Value = "2.24"
$ws.Range("E51").Value = "  +7.52%  "

# Restore default cell style on the Price column (content-only change, no formatting change intended)
$ws.Range("D2:D51").Style = "Normal"

